# WAT.xlsx - "Add WAT21 & WAT22 test script"
#
# Changes applied to the "Test Cases" sheet (sheet1 / ActiveSheet):
#   1. WAT02 (row 3) description text corrected.
#   2. WAT04 (row 5) JIRA ID list corrected.
#   3. WAT18 (row 19) description text corrected (dropped the stray
#      " \n * " typo) and its row height reverts to the default.
#   4. Two brand-new test rows are appended: WAT21 (row 22) and
#      WAT22 (row 23).
#   5. Column C is widened to fit the new, longer text.
#   6. Selection is left on C28:C29, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. WAT02 - fix the description copy.
$ws.Range("C3").Value = "Verify that user is able to search for an Author cluster u7sing both Last name and First  name and navigate to author search result page."

# 2. WAT04 - correct the linked JIRA ids.
$ws.Range("B5").Value = "WAT-156||WAT-160||WAT-195"

# 3. WAT18 - remove the embedded newline/typo from the description,
#    then let Excel recompute the row height (drops the custom ht="30").
$ws.Range("C19").Value = 'Verify that "Add alternative name" button should be disabled (Gryed out) until a single letter is entered in last name field'
$ws.Rows.Item(19).AutoFit()

# 4. New rows: WAT21 & WAT22.
$ws.Range("A22").Value = "WAT21"
$ws.Range("B22").Value = "WAT-190"
$ws.Range("C22").Value = "Verify that when there is no result for the search name entered by User then correct error message should be displayed - Last name"
$ws.Range("D22").Value = "Y"

$ws.Range("A23").Value = "WAT22"
$ws.Range("B23").Value = "WAT-190"
$ws.Range("C23").Value = "Verify that when there is no result for the search name entered by User then correct error message should be displayed - First name"
$ws.Range("D23").Value = "Y"

# Match the bordered look of the rest of the table for the new rows.
$ws.Range("A22:E23").Borders.LineStyle = 1

# 5. Column C needs to be wide enough for the longest entry now.
$ws.Columns.Item(3).ColumnWidth = 138.140625

# 6. Leave the selection where the author left it.
$ws.Range("C28:C29").Select()
